$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to Text format so numeric-looking price strings
# (e.g. "21.89") are preserved as text, matching the source data which
# stores every Price/Volume cell as a string.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "30.388.34"
$ws.Range("D3").Value = "1.872.18"
$ws.Range("D4").Value = "0.9995"
$ws.Range("D5").Value = "244.24"
$ws.Range("D7").Value = "0.4704"
$ws.Range("D8").Value = "0.2873"
$ws.Range("D9").Value = "0.06491"
$ws.Range("D10").Value = "21.89"
$ws.Range("D11").Value = "99.96"
$ws.Range("D12").Value = "0.07790"
$ws.Range("D13").Value = "1.872.68"
$ws.Range("D14").Value = "0.7301"
$ws.Range("D15").Value = "5.174"
$ws.Range("D16").Value = "286.55"
$ws.Range("D17").Value = "30.368.25"
$ws.Range("D19").Value = "0.9997"
$ws.Range("D20").Value = "0.000007487"
$ws.Range("D21").Value = "2.114.06"
$ws.Range("D22").Value = "5.299"
$ws.Range("D23").Value = "0.9992"
$ws.Range("D24").Value = "6.338"
$ws.Range("D25").Value = "163.32"
$ws.Range("D26").Value = "9.052"
$ws.Range("D27").Value = "19.01"
$ws.Range("D28").Value = "1.896"
$ws.Range("D29").Value = "0.09650"
$ws.Range("D30").Value = "1.318"
$ws.Range("D31").Value = "1.485"
$ws.Range("D32").Value = "4.234"
$ws.Range("D33").Value = "4.149"
$ws.Range("D35").Value = "1.126"
$ws.Range("D36").Value = "0.6885"
$ws.Range("D37").Value = "2.725"
$ws.Range("D40").Value = "76.07"
$ws.Range("D41").Value = "6.270"
$ws.Range("D42").Value = "1.957"
$ws.Range("D43").Value = "0.4228"
$ws.Range("D44").Value = "0.9989"
$ws.Range("D45").Value = "0.8235"
$ws.Range("D46").Value = "101.07"
$ws.Range("D47").Value = "9.794"
$ws.Range("D48").Value = "7.018"
$ws.Range("D49").Value = "35.00"
$ws.Range("D50").Value = "0.05766"
$ws.Range("D51").Value = "887.64"

$ws.Range("E2").Value = "  -1.11%  "
$ws.Range("E3").Value = "  -0.94%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("E5").Value = "  -2.12%  "
$ws.Range("E6").Value = "  -0.04%  "
$ws.Range("E7").Value = "  -1.24%  "
$ws.Range("E8").Value = "  -2.08%  "
$ws.Range("E9").Value = "  -0.69%  "
$ws.Range("E11").Value = "  +3.20%  "
$ws.Range("E12").Value = "  +0.41%  "
$ws.Range("E13").Value = "  -0.88%  "
$ws.Range("E14").Value = "  -1.33%  "
$ws.Range("E15").Value = "  -1.60%  "
$ws.Range("E16").Value = "  +3.30%  "
$ws.Range("E17").Value = "  -1.42%  "
$ws.Range("E18").Value = "  -0.61%  "
$ws.Range("E19").Value = "  -0.01%  "
$ws.Range("E20").Value = "  -1.20%  "
$ws.Range("E21").Value = "  -0.98%  "
$ws.Range("E22").Value = "  -0.70%  "
$ws.Range("E23").Value = "  -0.07%  "
$ws.Range("E24").Value = "  +1.43%  "
$ws.Range("E25").Value = "  -0.70%  "
$ws.Range("E26").Value = "  -2.18%  "
$ws.Range("E27").Value = "  +0.84%  "
$ws.Range("E28").Value = "  -1.73%  "
$ws.Range("E29").Value = "  -0.89%  "
$ws.Range("E30").Value = "  -2.10%  "
$ws.Range("E31").Value = "  -1.14%  "
$ws.Range("E32").Value = "  -1.81%  "
$ws.Range("E33").Value = "  -1.38%  "
$ws.Range("E34").Value = "  -1.54%  "
$ws.Range("E35").Value = "  -0.07%  "
$ws.Range("E36").Value = "  -1.72%  "
$ws.Range("E37").Value = "  +0.10%  "
$ws.Range("E38").Value = "  -1.17%  "
$ws.Range("E39").Value = "  +1.18%  "
$ws.Range("E40").Value = "  -0.59%  "
$ws.Range("E41").Value = "  -1.26%  "
$ws.Range("E42").Value = "  -3.71%  "
$ws.Range("E43").Value = "  -0.82%  "
$ws.Range("E44").Value = "  -0.09%  "
$ws.Range("E45").Value = "  -2.16%  "
$ws.Range("E46").Value = "  -0.96%  "
$ws.Range("E47").Value = "  +3.89%  "
$ws.Range("E48").Value = "  -0.67%  "
$ws.Range("E49").Value = "  -2.14%  "
$ws.Range("E50").Value = "  -0.06%  "
$ws.Range("E51").Value = "  -3.83%  "

# Reset the style back to Normal so no stray number-format style
# index is left attached to the cells (matches original formatting).
$ws.Range("D2:D51").Style = "Normal"
